$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author regraded four students from "B" to "C".
$ws.Range("B18").Value = "C"
$ws.Range("B75").Value = "C"
$ws.Range("B95").Value = "C"
$ws.Range("B100").Value = "C"

# Restore the view state (scroll position / active selection) that was
# captured at save time: the window had scrolled so row 22 is the top
# visible row, and F36 was the active/selected cell.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F36").Select()
